# Zakat & Tax SignUp Developed.
# Updates the Controller run-status cell and expands the ZakatAndTax
# "Consortium Company" sign-up data sheet with the new AccountType /
# Consortium company fields.

$wb = $excel.ActiveWorkbook

# --- Controller sheet: mark the test run as Passed -------------------------
$controller = $wb.Worksheets.Item("Controller")
$controller.Range("H3").Value = "Passed"

# --- ZakatAndTax sheet: add the Consortium-company sign-up columns ---------
$ws = $wb.Worksheets.Item("ZakatAndTax")

# Existing "OtherAccountCategory" cell (E2) was blank - now holds the
# account-category value used for the Consortium flow.
$ws.Range("E2").Value = "Consortium"

# New headers (row 1) - copy the header formatting from D1 (same row style
# used for every header cell) onto the new header cells first
$ws.Range("D1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "Company_ID"
$ws.Range("G1").Value = "Company_Name"
$ws.Range("H1").Value = "Contract_Number"
$ws.Range("I1").Value = "Mobile_Number"
$ws.Range("J1").Value = "Email"
$ws.Range("K1").Value = "Confirm_Email"

# New data (row 2) - Company_ID must stay text (leading-safe, 11 digits)
$ws.Range("F2").Value = "'71234567890"
$ws.Range("G2").Value = "TestCrewConsortium"
$ws.Range("H2").Value = 123456
$ws.Range("I2").Value = 555000001
$ws.Range("J2").Value = "TestCrewConsortium@lab.testcrew.com"
$ws.Range("K2").Value = "TestCrewConsortium@lab.testcrew.com"

# Email / Confirm_Email get mailto hyperlinks, like the existing URL column
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:TestCrewConsortium@lab.testcrew.com")
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:TestCrewConsortium@lab.testcrew.com")

# Reflect the new columns in the view: scroll right and select F3
$ws.Activate() | Out-Null
$ws.Range("D1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F3").Select() | Out-Null
